$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.09286699999999999
$ws.Range("H2").Value = 0.278601
$ws.Range("I2").Value = 0.003009076821730935
$ws.Range("J2").Value = 0.003071957783644885
$ws.Range("M2").Value = 28.85518433333334
$ws.Range("N2").Value = 86.56555300000001
$ws.Range("O2").Value = 0.1999651185353207
$ws.Range("P2").Value = 0.2044513327926365
$ws.Range("Q2").Value = 2.679694403483667
$ws.Range("R2").Value = 24.117249631353
$ws.Range("S2").Value = 0.0006017104033393123
$ws.Range("T2").Value = 0.0006280658631489104

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.09286699999999999
$ws.Range("H3").Value = 0.278601
$ws.Range("I3").Value = 0.003009076821730935
$ws.Range("J3").Value = 0.003071957783644885
$ws.Range("O3").Value = 0.3546352265743414
$ws.Range("P3").Value = 0.3625914622481308
$ws.Range("Q3").Value = 4.752399012838667
$ws.Range("R3").Value = 42.771591115548
$ws.Range("S3").Value = 0.001067124640454149
$ws.Range("T3").Value = 0.001113865664736326

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.09286699999999999
$ws.Range("H4").Value = 0.278601
$ws.Range("I4").Value = 0.003009076821730935
$ws.Range("J4").Value = 0.003071957783644885
$ws.Range("M4").Value = 29.393479
$ws.Range("N4").Value = 88.180437
$ws.Range("O4").Value = 0.2036954761578358
$ws.Range("P4").Value = 0.2082653809291453
$ws.Range("Q4").Value = 2.729684214293
$ws.Range("R4").Value = 24.567157928637
$ws.Range("S4").Value = 0.0006129353359979901
$ws.Range("T4").Value = 0.0006397824580090548

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.09286699999999999
$ws.Range("H5").Value = 0.278601
$ws.Range("I5").Value = 0.003009076821730935
$ws.Range("J5").Value = 0.003071957783644885
$ws.Range("M5").Value = 9.499066500000001
$ws.Range("N5").Value = 18.998133
$ws.Range("O5").Value = 0.0658280999596015
$ws.Range("P5").Value = 0.04486996822421697
$ws.Range("Q5").Value = 0.8821498086555001
$ws.Range("R5").Value = 5.292898851933001
$ws.Range("S5").Value = 0.0001980818098070239
$ws.Range("T5").Value = 0.000137838648138282

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.09286699999999999
$ws.Range("H6").Value = 0.278601
$ws.Range("I6").Value = 0.003009076821730935
$ws.Range("J6").Value = 0.003071957783644885
$ws.Range("M6").Value = 25.37910966666666
$ws.Range("N6").Value = 76.137329
$ws.Range("O6").Value = 0.1758760787729007
$ws.Range("P6").Value = 0.1798218558058706
$ws.Range("Q6").Value = 2.356881777414333
$ws.Range("R6").Value = 21.211935996729
$ws.Range("S6").Value = 0.0005292246321324595
$ws.Range("T6").Value = 0.0005524051496123125

# Row 7
$ws.Range("I7").Value = 0.9272539658256183
$ws.Range("J7").Value = 0.9466308793322996
$ws.Range("M7").Value = 28.85518433333334
$ws.Range("N7").Value = 86.56555300000001
$ws.Range("O7").Value = 0.1999651185353207
$ws.Range("P7").Value = 0.2044513327926365
$ws.Range("Q7").Value = 825.7540136185748
$ws.Range("R7").Value = 7431.786122567173
$ws.Range("S7").Value = 0.1854184491886659
$ws.Range("T7").Value = 0.1935399449421541

# Row 8
$ws.Range("I8").Value = 0.9272539658256183
$ws.Range("J8").Value = 0.9466308793322996
$ws.Range("O8").Value = 0.3546352265743414
$ws.Range("P8").Value = 0.3625914622481308
$ws.Range("S8").Value = 0.3288369202625248
$ws.Range("T8").Value = 0.3432402747463324

# Row 9
$ws.Range("I9").Value = 0.9272539658256183
$ws.Range("J9").Value = 0.9466308793322996
$ws.Range("M9").Value = 29.393479
$ws.Range("N9").Value = 88.180437
$ws.Range("O9").Value = 0.2036954761578358
$ws.Range("P9").Value = 0.2082653809291453
$ws.Range("Q9").Value = 841.1584891670464
$ws.Range("R9").Value = 7570.426402503416
$ws.Range("S9").Value = 0.188877438088091
$ws.Range("T9").Value = 0.1971504406834331

# Row 10
$ws.Range("I10").Value = 0.9272539658256183
$ws.Range("J10").Value = 0.9466308793322996
$ws.Range("M10").Value = 9.499066500000001
$ws.Range("N10").Value = 18.998133
$ws.Range("O10").Value = 0.0658280999596015
$ws.Range("P10").Value = 0.04486996822421697
$ws.Range("Q10").Value = 271.8364990288255
$ws.Range("R10").Value = 1631.018994172953
$ws.Range("S10").Value = 0.06103936675030571
$ws.Range("T10").Value = 0.04247529747570285

# Row 11
$ws.Range("I11").Value = 0.9272539658256183
$ws.Range("J11").Value = 0.9466308793322996
$ws.Range("M11").Value = 25.37910966666666
$ws.Range("N11").Value = 76.137329
$ws.Range("O11").Value = 0.1758760787729007
$ws.Range("P11").Value = 0.1798218558058706
$ws.Range("Q11").Value = 726.2785580304432
$ws.Range("R11").Value = 6536.507022273989
$ws.Range("S11").Value = 0.163081791536031
$ws.Range("T11").Value = 0.1702249214846773

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.243062
$ws.Range("H12").Value = 0.729186
$ws.Range("I12").Value = 0.007875695677081898
$ws.Range("J12").Value = 0.00804027483183793
$ws.Range("M12").Value = 28.85518433333334
$ws.Range("N12").Value = 86.56555300000001
$ws.Range("O12").Value = 0.1999651185353207
$ws.Range("P12").Value = 0.2044513327926365
$ws.Range("Q12").Value = 7.013598814428668
$ws.Range("R12").Value = 63.12238932985801
$ws.Range("S12").Value = 0.001574864419615794
$ws.Range("T12").Value = 0.001643844905388356

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.243062
$ws.Range("H13").Value = 0.729186
$ws.Range("I13").Value = 0.007875695677081898
$ws.Range("J13").Value = 0.00804027483183793
$ws.Range("O13").Value = 0.3546352265743414
$ws.Range("P13").Value = 0.3625914622481308
$ws.Range("Q13").Value = 12.43851539145867
$ws.Range("R13").Value = 111.946638523128
$ws.Range("S13").Value = 0.0027929991208725
$ws.Range("T13").Value = 0.002915335008152959

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.243062
$ws.Range("H14").Value = 0.729186
$ws.Range("I14").Value = 0.007875695677081898
$ws.Range("J14").Value = 0.00804027483183793
$ws.Range("M14").Value = 29.393479
$ws.Range("N14").Value = 88.180437
$ws.Range("O14").Value = 0.2036954761578358
$ws.Range("P14").Value = 0.2082653809291453
$ws.Range("Q14").Value = 7.144437792698
$ws.Range("R14").Value = 64.299940134282
$ws.Range("S14").Value = 0.001604243581017406
$ws.Range("T14").Value = 0.001674510900627746

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.243062
$ws.Range("H15").Value = 0.729186
$ws.Range("I15").Value = 0.007875695677081898
$ws.Range("J15").Value = 0.00804027483183793
$ws.Range("M15").Value = 9.499066500000001
$ws.Range("N15").Value = 18.998133
$ws.Range("O15").Value = 0.0658280999596015
$ws.Range("P15").Value = 0.04486996822421697
$ws.Range("Q15").Value = 2.308862101623
$ws.Range("R15").Value = 13.853172609738
$ws.Range("S15").Value = 0.0005184420822823485
$ws.Range("T15").Value = 0.0003607668762185393

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.243062
$ws.Range("H16").Value = 0.729186
$ws.Range("I16").Value = 0.007875695677081898
$ws.Range("J16").Value = 0.00804027483183793
$ws.Range("M16").Value = 25.37910966666666
$ws.Range("N16").Value = 76.137329
$ws.Range("O16").Value = 0.1758760787729007
$ws.Range("P16").Value = 0.1798218558058706
$ws.Range("Q16").Value = 6.168697153799333
$ws.Range("R16").Value = 55.51827438419399
$ws.Range("S16").Value = 0.001385146473293849
$ws.Range("T16").Value = 0.001445817141450331

# Row 17
$ws.Range("G17").Value = 1.8951925
$ws.Range("H17").Value = 3.790385
$ws.Range("I17").Value = 0.06140803325689756
$ws.Range("J17").Value = 0.04179418847656979
$ws.Range("M17").Value = 28.85518433333334
$ws.Range("N17").Value = 86.56555300000001
$ws.Range("O17").Value = 0.1999651185353207
$ws.Range("P17").Value = 0.2044513327926365
$ws.Range("Q17").Value = 54.68612893465085
$ws.Range("R17").Value = 328.116773607905
$ws.Range("S17").Value = 0.01227946464923643
$ws.Range("T17").Value = 0.008544877537021343

# Row 18
$ws.Range("G18").Value = 1.8951925
$ws.Range("H18").Value = 3.790385
$ws.Range("I18").Value = 0.06140803325689756
$ws.Range("J18").Value = 0.04179418847656979
$ws.Range("O18").Value = 0.3546352265743414
$ws.Range("P18").Value = 0.3625914622481308
$ws.Range("Q18").Value = 96.98505352966335
$ws.Range("R18").Value = 581.91032117798
$ws.Range("S18").Value = 0.02177745178754456
$ws.Range("T18").Value = 0.01515421591319342

# Row 19
$ws.Range("G19").Value = 1.8951925
$ws.Range("H19").Value = 3.790385
$ws.Range("I19").Value = 0.06140803325689756
$ws.Range("J19").Value = 0.04179418847656979
$ws.Range("M19").Value = 29.393479
$ws.Range("N19").Value = 88.180437
$ws.Range("O19").Value = 0.2036954761578358
$ws.Range("P19").Value = 0.2082653809291453
$ws.Range("Q19").Value = 55.7063009497075
$ws.Range("R19").Value = 334.237805698245
$ws.Range("S19").Value = 0.01250853857417997
$ws.Range("T19").Value = 0.0087042825836973

# Row 20
$ws.Range("G20").Value = 1.8951925
$ws.Range("H20").Value = 3.790385
$ws.Range("I20").Value = 0.06140803325689756
$ws.Range("J20").Value = 0.04179418847656979
$ws.Range("M20").Value = 9.499066500000001
$ws.Range("N20").Value = 18.998133
$ws.Range("O20").Value = 0.0658280999596015
$ws.Range("P20").Value = 0.04486996822421697
$ws.Range("Q20").Value = 18.00255958780125
$ws.Range("R20").Value = 72.01023835120502
$ws.Range("S20").Value = 0.004042374151557586
$ws.Range("T20").Value = 0.001875303908900621

# Row 21
$ws.Range("G21").Value = 1.8951925
$ws.Range("H21").Value = 3.790385
$ws.Range("I21").Value = 0.06140803325689756
$ws.Range("J21").Value = 0.04179418847656979
$ws.Range("M21").Value = 25.37910966666666
$ws.Range("N21").Value = 76.137329
$ws.Range("O21").Value = 0.1758760787729007
$ws.Range("P21").Value = 0.1798218558058706
$ws.Range("Q21").Value = 48.09829829694417
$ws.Range("R21").Value = 288.589789781665
$ws.Range("S21").Value = 0.01080020409437902
$ws.Range("T21").Value = 0.007515508533757112

# Row 22
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.3333333333333333
$ws.Range("G22").Value = 0.01398766666666667
$ws.Range("H22").Value = 0.041963
$ws.Range("I22").Value = 0.000453228418671488
$ws.Range("J22").Value = 0.0004626995756479349
$ws.Range("M22").Value = 28.85518433333334
$ws.Range("N22").Value = 86.56555300000001
$ws.Range("O22").Value = 0.1999651185353207
$ws.Range("P22").Value = 0.2044513327926365
$ws.Range("Q22").Value = 0.403616700059889
$ws.Range("R22").Value = 3.632550300539001
$ws.Range("S22").Value = 0.00009062987446322005
$ws.Range("T22").Value = 0.00009459954492380764

# Row 23
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.3333333333333333
$ws.Range("G23").Value = 0.01398766666666667
$ws.Range("H23").Value = 0.041963
$ws.Range("I23").Value = 0.000453228418671488
$ws.Range("J23").Value = 0.0004626995756479349
$ws.Range("O23").Value = 0.3546352265743414
$ws.Range("P23").Value = 0.3625914622481308
$ws.Range("Q23").Value = 0.7158083415915556
$ws.Range("R23").Value = 6.442275074324001
$ws.Range("S23").Value = 0.0001607307629454936
$ws.Range("T23").Value = 0.0001677709157157743

# Row 24
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.01398766666666667
$ws.Range("H24").Value = 0.041963
$ws.Range("I24").Value = 0.000453228418671488
$ws.Range("J24").Value = 0.0004626995756479349
$ws.Range("M24").Value = 29.393479
$ws.Range("N24").Value = 88.180437
$ws.Range("O24").Value = 0.2036954761578358
$ws.Range("P24").Value = 0.2082653809291453
$ws.Range("Q24").Value = 0.4111461864256667
$ws.Range("R24").Value = 3.700315677831
$ws.Range("S24").Value = 0.00009232057854955173
$ws.Range("T24").Value = 0.00009636430337807104

# Row 25
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0.3333333333333333
$ws.Range("G25").Value = 0.01398766666666667
$ws.Range("H25").Value = 0.041963
$ws.Range("I25").Value = 0.000453228418671488
$ws.Range("J25").Value = 0.0004626995756479349
$ws.Range("M25").Value = 9.499066500000001
$ws.Range("N25").Value = 18.998133
$ws.Range("O25").Value = 0.0658280999596015
$ws.Range("P25").Value = 0.04486996822421697
$ws.Range("Q25").Value = 0.1328697758465
$ws.Range("R25").Value = 0.7972186550790001
$ws.Range("S25").Value = 0.00002983516564883883
$ws.Range("T25").Value = 0.00002076131525668152

# Row 26
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0.3333333333333333
$ws.Range("G26").Value = 0.01398766666666667
$ws.Range("H26").Value = 0.041963
$ws.Range("I26").Value = 0.000453228418671488
$ws.Range("J26").Value = 0.0004626995756479349
$ws.Range("M26").Value = 25.37910966666666
$ws.Range("N26").Value = 76.137329
$ws.Range("O26").Value = 0.1758760787729007
$ws.Range("P26").Value = 0.1798218558058706
$ws.Range("Q26").Value = 0.3549945263141111
$ws.Range("R26").Value = 3.194950736827
$ws.Range("S26").Value = 0.00007971203706438383
$ws.Range("T26").Value = 0.00008320349637360048
